$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.665.26"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.438.32"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "2.890.47"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "62.550.07"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "2.443.24"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.85%  "
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "578.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("E27").Value = "  +8.77%  "
$ws.Range("D28").Value = "2.557.85"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.146"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.382"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "148.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  +9.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "148.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0534"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.601"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0919"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +4.76%  "
